$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the byte order in row 1 (B1:I1): was 7,6,5,4,3,2,1,0 -> now 0,1,2,3,4,5,6,7
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 6
$ws.Range("I1").Value = 7

# Update selection to match new active cell
$ws.Range("B4").Select()

# Update the window position
$excel.ActiveWindow.WindowState = -4143
$wb.Windows.Item(1).Left = 1950
$wb.Windows.Item(1).Top = 1950
